# Automated "Disponibilidad" checker run: appends one fresh batch of 14
# monitor results (rows 1052-1065) below the last existing batch
# (which ended at row 1051), and nudges the D-column timestamp of the
# previous batch (rows 1038-1051) by the tiny recalculation delta that
# the authoring tool produced on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Re-stamp the previous batch's "Ultimo" timestamp (rows 1038-1051) ---
$prevTimestamp = 44233.23612936342
for ($r = 1038; $r -le 1051; $r++) {
    $ws.Cells.Item($r, 4).Value = $prevTimestamp
}

# --- 2. Append the new batch (rows 1052-1065) ---
$newTimestamp = 44233.25741231485

$batch = @(
    @{Row=1052; Name="Odoo";               Url="https://www.dataintelligence-group.com/"},
    @{Row=1053; Name="Blackbox";            Url="https://serviciodashboard.azurewebsites.net/"},
    @{Row=1054; Name="PowerBI";             Url="https://powerbi.microsoft.com/es-es/"},
    @{Row=1055; Name="Dropbox";             Url="https://www.dropbox.com/"},
    @{Row=1056; Name="Odoo";               Url="https://dataintelligence.store/"},
    @{Row=1057; Name="GEE";                 Url="https://app-data-i.users.earthengine.app/"},
    @{Row=1058; Name="UtilidadesOdoo";      Url="https://odooutil.azurewebsites.net/"},
    @{Row=1059; Name="Filtros Dashboard";   Url="https://filtradordashboard.azurewebsites.net/"},
    @{Row=1060; Name="MapStore";            Url="https://ide.dataintelligence-group.com/mapstore/#/"},
    @{Row=1061; Name="GeoServer";           Url="https://ide.dataintelligence-group.com/geoserver/web/?0"},
    @{Row=1062; Name="Tomcat";              Url="https://ide.dataintelligence-group.com/"},
    @{Row=1063; Name="Shiny";               Url="https://rpubs.com/dataintelligence/"},
    @{Row=1064; Name="Github";              Url="https://github.com/Sud-Austral/"},
    @{Row=1065; Name="EZ Exporter";         Url="https://ezexporter.highviewapps.com/exports/export-profile/"}
)

foreach ($entry in $batch) {
    $r = $entry.Row

    $ws.Cells.Item($r, 1).Value = $entry.Name

    $linkCell = $ws.Cells.Item($r, 2)
    $linkCell.Value = $entry.Url
    $ws.Hyperlinks.Add($linkCell, $entry.Url) | Out-Null
    $linkCell.Style = "Hyperlink"

    $ws.Cells.Item($r, 3).Value = "Disponible"

    $dateCell = $ws.Cells.Item($r, 4)
    $dateCell.Value = $newTimestamp
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

Write-Host "Actualizacion completa: filas 1052-1065 agregadas."
